# atualizacao de dados: inclui os dias 23, 24 e 25 de maio/2025 que faltavam
# na planilha de faturamento diario. As linhas novas sao inseridas logo apos
# o ultimo dia de maio ja existente (linha 23), empurrando abril/marco/fevereiro
# para baixo em 3 linhas (sem alterar nenhum desses valores).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insere 3 novas linhas em branco a partir da linha 24 (antes do primeiro dia de abril)
$ws.Range("A24:E26").EntireRow.Insert()

# Novos registros de maio/2025 (dias 23-25)
$novosDados = @(
    @(23, 23820.29, 5, 2025, "05/2025"),
    @(24, 34922.72, 5, 2025, "05/2025"),
    @(25, 1686.9,   5, 2025, "05/2025")
)

$r = 24
foreach ($linha in $novosDados) {
    $ws.Cells.Item($r, 1).Value = $linha[0]
    $ws.Cells.Item($r, 2).Value = $linha[1]
    $ws.Cells.Item($r, 3).Value = $linha[2]
    $ws.Cells.Item($r, 4).Value = $linha[3]
    $ws.Cells.Item($r, 5).Value = $linha[4]
    $r++
}
